$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.525.27"
$ws.Range("E2").Value = "  +5.28%  "
$ws.Range("D3").Value = "3.632.92"
$ws.Range("E3").Value = "  +5.55%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.87"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "191.41"
$ws.Range("E6").Value = "  +3.61%  "
$ws.Range("E7").Value = "  +1.98%  "
$ws.Range("D8").Value = "3.628.71"
$ws.Range("E8").Value = "  +5.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  +2.61%  "
$ws.Range("E11").Value = "  +2.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.33"
$ws.Range("E12").Value = "  +3.46%  "
$ws.Range("E13").Value = "  +3.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.91"
$ws.Range("E14").Value = "  +5.09%  "
$ws.Range("D15").Value = "4.203.66"
$ws.Range("E15").Value = "  +5.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.75"
$ws.Range("E16").Value = "  +5.78%  "
$ws.Range("D17").Value = "3.621.44"
$ws.Range("E17").Value = "  +5.51%  "
$ws.Range("D18").Value = "70.420.83"
$ws.Range("E18").Value = "  +5.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.68"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("E21").Value = "  +4.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "488.50"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("E23").Value = "  +15.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.34"
$ws.Range("E24").Value = "  +4.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.45"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.03"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("E27").Value = "  +6.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.31"
$ws.Range("E28").Value = "  +2.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.70"
$ws.Range("E29").Value = "  +5.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.11"
$ws.Range("E30").Value = "  +5.62%  "
$ws.Range("E31").Value = "  +8.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "635.39"
$ws.Range("E32").Value = "  +7.63%  "
$ws.Range("E33").Value = "  +5.13%  "
$ws.Range("E34").Value = "  +7.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "66.24"
$ws.Range("E35").Value = "  +3.37%  "
$ws.Range("B36").Value = "TheGraph"
$ws.Range("C36").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.413"
$ws.Range("E36").Value = "  +7.00%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.79"
$ws.Range("E37").Value = "  +6.38%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0822"
$ws.Range("E38").Value = "  +6.87%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "3.309.29"
$ws.Range("E42").Value = "  +3.70%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.82"
$ws.Range("E43").Value = "  +10.97%  "
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.10"
$ws.Range("E44").Value = "  +6.21%  "
$ws.Range("E45").Value = "  +4.93%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.139"
$ws.Range("E46").Value = "  +2.84%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.29"
$ws.Range("E47").Value = "  +2.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.11"
$ws.Range("E48").Value = "  +4.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.72"
$ws.Range("E49").Value = "  -2.11%  "
$ws.Range("E50").Value = "  +3.70%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  +0.02%  "
